# Apply the "grouping and column split" edit:
# - Split Significant_Diseases (col N) into 3 new columns AD/AE/AF
#   (Significant_Diseases_part_1/2/3), splitting on ",".
#   Special case: if the first part trims to "None", blank it out.
# - Split Friends_To_Talk (col I) into 2 new columns AG/AH
#   (Friends_To_Talk_part_1/2), splitting on ",".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
# Copy the existing header formatting (bold, border, centered) from A1
# onto the 5 new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AH1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value2 = "Significant_Diseases_part_1"
$ws.Range("AE1").Value2 = "Significant_Diseases_part_2"
$ws.Range("AF1").Value2 = "Significant_Diseases_part_3"
$ws.Range("AG1").Value2 = "Friends_To_Talk_part_1"
$ws.Range("AH1").Value2 = "Friends_To_Talk_part_2"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {

    $diseases = $ws.Cells.Item($r, 14).Value2   # column N = 14
    $friends  = $ws.Cells.Item($r, 9).Value2    # column I = 9

    $dParts = @("", "", "")
    if ($diseases -ne $null -and $diseases -ne "") {
        $split = $diseases.Split(",")
        for ($i = 0; $i -lt $split.Length -and $i -lt 3; $i++) {
            $dParts[$i] = $split[$i]
        }
        if ($dParts[0].Trim() -eq "None") {
            $dParts[0] = ""
        }
    }

    $fParts = @("", "")
    if ($friends -ne $null -and $friends -ne "") {
        $split2 = $friends.Split(",")
        for ($i = 0; $i -lt $split2.Length -and $i -lt 2; $i++) {
            $fParts[$i] = $split2[$i]
        }
    }

    $ws.Cells.Item($r, 30).Value2 = $dParts[0]   # AD
    $ws.Cells.Item($r, 31).Value2 = $dParts[1]   # AE
    $ws.Cells.Item($r, 32).Value2 = $dParts[2]   # AF
    $ws.Cells.Item($r, 33).Value2 = $fParts[0]   # AG
    $ws.Cells.Item($r, 34).Value2 = $fParts[1]   # AH
}

Write-Output "done"
